$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Outbound flight (row 6)
$ws.Range("A6").Value = "Frontier"
$ws.Range("B6").Value = "6:05 AM – 8:04 AM"
$ws.Range("D6").Value = "2 hr 59 min"
$ws.Range("E6").Value = 126

# Return flight (row 7)
$ws.Range("A7").Value = "Spirit"
$ws.Range("B7").Value = "11:10 AM – 1:03 PM"
$ws.Range("D7").Value = "2 hr 53 min"
$ws.Range("E7").Value = 191
